$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the interest rate input on the amortisation schedule
# (was 7.5%, now 4.2%) - all dependent schedule/summary formulas
# recalculate automatically.
$ws.Range("B4").Value = 4.2

# Reflect the author's new active cell selection
$ws.Range("B6").Select()
